$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily price-refresh update (GitHub Actions cron job).
#
# All Price (D) / Volume(1h) (E) cells in this sheet are stored as literal TEXT
# (t="inlineStr"), not numbers -- e.g. "26.353.74" or "19.51" are display strings,
# not numeric values. Assigning a plain numeric-looking string straight to
# Range.Value makes Excel auto-coerce it into a real number (losing the exact
# text and switching the cell type), so for any new Price value that parses as a
# plain number we force the cell to Text format first, then restore the original
# (default/unstyled) cell style afterwards so no stray formatting is left behind.
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $originalStyle
}

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.353.74"
$ws.Range("E2").Value = "  -1.10%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.589.53"
$ws.Range("E3").Value = "  -0.54%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.43%  "

# Row 5: BNB
Set-TextValue "D5" "210.13"
$ws.Range("E5").Value = "  -0.29%  "

# Row 6: XRP
$ws.Range("E6").Value = "  -1.16%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.40%  "

# Row 8: Dogecoin
$ws.Range("E8").Value = "  -0.84%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  -0.51%  "

# Row 10: Solana
Set-TextValue "D10" "19.51"
$ws.Range("E10").Value = "  -0.39%  "

# Row 11: TRON
Set-TextValue "D11" "0.0845"
$ws.Range("E11").Value = "  +0.01%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.814.14"
$ws.Range("E12").Value = "  -0.46%  "

# Row 13: WrappedEther
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.07"
$ws.Range("E13").Value = "  +0.64%  "

# Row 14: Polkadot
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.564.09"
$ws.Range("E14").Value = "  -2.18%  "

# Row 15: Polygon
$ws.Range("E15").Value = "  -0.92%  "

# Row 16: Litecoin
Set-TextValue "D16" "64.33"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "26.361.58"
$ws.Range("E17").Value = "  -0.96%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -1.63%  "

# Row 19: Chainlink
Set-TextValue "D19" "7.45"
$ws.Range("E19").Value = "  +4.99%  "

# Row 20: BitcoinCash
Set-TextValue "D20" "210.54"
$ws.Range("E20").Value = "  +1.07%  "

# Row 21: Dai
$ws.Range("E21").Value = "  -0.45%  "

# Row 22: Uniswap
$ws.Range("E22").Value = "  -0.17%  "

# Row 23: Toncoin
Set-TextValue "D23" "2.14"
$ws.Range("E23").Value = "  -3.96%  "

# Row 24: Avalanche
Set-TextValue "D24" "8.91"
$ws.Range("E24").Value = "  -0.28%  "

# Row 25: Monero
Set-TextValue "D25" "145.01"
$ws.Range("E25").Value = "  +0.92%  "

# Row 26: BinanceUSD
$ws.Range("E26").Value = "  -0.43%  "

# Row 27: Cosmos
$ws.Range("E27").Value = "  -1.26%  "

# Row 29: EthereumClassic
Set-TextValue "D29" "15.24"

# Row 30: Hedera
$ws.Range("E30").Value = "  -0.38%  "

# Row 31: PancakeSwap
$ws.Range("E31").Value = "  -0.22%  "

# Row 32: Filecoin
$ws.Range("E32").Value = "  -0.78%  "

# Row 34: Maker
$ws.Range("D34").Value = "1.306.03"
$ws.Range("E34").Value = "  +2.14%  "

# Row 35: ImmutableX
Set-TextValue "D35" "0.613"
$ws.Range("E35").Value = "  +2.55%  "

# Row 36: HuobiToken
$ws.Range("E36").Value = "  -1.70%  "

# Row 37: LidoDAOToken
Set-TextValue "D37" "1.47"
$ws.Range("E37").Value = "  -0.90%  "

# Row 38: VeChain
$ws.Range("E38").Value = "  +0.41%  "

# Row 39: WEMIXToken
$ws.Range("E39").Value = "  -13.99%  "

# Row 40: ARBITRUM
$ws.Range("E40").Value = "  -1.63%  "

# Row 41: PaxDollar
$ws.Range("E41").Value = "  -0.39%  "

# Row 42: FraxShare
$ws.Range("E42").Value = "  +3.68%  "

# Row 43: MXToken
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "0.767"
$ws.Range("E43").Value = "  -0.98%  "

# Row 44: TrustWalletToken
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D44" "2.13"
$ws.Range("E44").Value = "  -1.30%  "

# Row 45: Aave
Set-TextValue "D45" "62.54"
$ws.Range("E45").Value = "  -0.16%  "

# Row 46: RocketPoolETH
$ws.Range("D46").Value = "1.725.67"
$ws.Range("E46").Value = "  -0.45%  "

# Row 47: Quant
Set-TextValue "D47" "87.88"
$ws.Range("E47").Value = "  -1.92%  "

# Row 48: RenderToken
Set-TextValue "D48" "1.49"
$ws.Range("E48").Value = "  -4.86%  "

# Row 49: Algorand
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.0505"
$ws.Range("E49").Value = "  -1.46%  "

# Row 50: Cronos
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0981"
$ws.Range("E50").Value = "  -4.52%  "

# Row 51: USDD
$ws.Range("E51").Value = "  -0.39%  "
